# Daily attendance processing - reorder "Recorded By" (column G) entries.
#
# The G column holds a comma-separated list of recorder identities, e.g.
#   "dnasr281@gmail.com, System"
#   "backup@backdoor.com, system, System"
#   "dnasr281@gmail.com, admin@admin.com"
#
# Processing moves administrative/system accounts ahead of the human
# recorder's email address:
#   - if "admin@admin.com" is present, it is moved to the front
#   - else if any case-insensitive "system" entry is present, all such
#     entries are moved to the front (preserving their relative order)
#   - otherwise (single entry, or no recognized system-like account) the
#     value is left unchanged

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

# Column G = "Recorded By" (column index 7)
$col = 7

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $raw = $cell.Value2

    if ($null -eq $raw) { continue }
    if ($raw -eq "") { continue }
    if ($raw -eq "Recorded By") { continue }

    $rawParts = $raw.Split(",")
    if ($rawParts.Length -lt 2) { continue }

    $parts = @()
    foreach ($p in $rawParts) {
        $parts += $p.Trim()
    }

    $newOrder = $null

    if ($parts -contains "admin@admin.com") {
        $rest = @()
        foreach ($p in $parts) {
            if ($p -ne "admin@admin.com") { $rest += $p }
        }
        $newOrder = @("admin@admin.com") + $rest
    }
    else {
        $sysLike = @()
        $other = @()
        foreach ($p in $parts) {
            if ($p.ToLower() -eq "system") { $sysLike += $p } else { $other += $p }
        }
        if ($sysLike.Length -gt 0) {
            $newOrder = $sysLike + $other
        }
    }

    if ($null -ne $newOrder) {
        $newVal = [string]::Join(", ", $newOrder)
        if ($newVal -ne $raw) {
            $cell.Value = $newVal
        }
    }
}
